# feat: add 2022-Q1 data
#
# The workbook has quarterly fund-holding detail sheets (2021-Q1 .. 2021-Q4)
# plus a rollup sheet "总计". This change:
#   1. Adds a new quarterly detail sheet "2022-Q1" (same layout as the other
#      quarters) positioned right before "总计".
#   2. Updates "总计" with a new first data row for 2022-Q1, pushing the
#      older quarters down by one row.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# --- Step 1: clone "总计" (placed right after it) - this clone will become
# the new, updated "总计" sheet. Doing this first lets the *original* "总计"
# sheet (and its sheetId/rId) become the new "2022-Q1" sheet, matching how
# Excel keeps ids stable on rename.
$total.Copy($null, $total)
$newTotal = $wb.Worksheets.Item("总计 (2)")

# --- Step 2: turn the original "总计" sheet into the new "2022-Q1" sheet.
$total.Cells.Clear()
$total.Name = "2022-Q1"

# Bring over the same layout/styles used by the other quarterly sheets.
$q4.Range("A1:H3").Copy($total.Range("A1"))

# Fill in this quarter's numbers (fund code/name & rank carry over unchanged;
# size/position figures are quarter-specific). Values are entered as text,
# matching how the source data is stored in the other quarter sheets.
$total.Range("D2").Value = "'0.20"
$total.Range("E2").Value = "'93.65"
$total.Range("F2").Value = "'3.47"
$total.Range("G2").Value = "'0.0069"

$total.Range("D3").Value = "'0.06"
$total.Range("E3").Value = "'93.65"
$total.Range("F3").Value = "'3.47"
$total.Range("G3").Value = "'0.0021"

# --- Step 3: rename the clone to "总计" and refresh its summary table with
# the new 2022-Q1 row, shifting the older quarters down by one row.
$newTotal.Name = "总计"

$newTotal.Range("A5:D5").Copy($newTotal.Range("A6"))

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 2
$newTotal.Range("D2").Value = 0.01

$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 2
$newTotal.Range("D3").Value = 0.01

$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 2
$newTotal.Range("D4").Value = 0.01

$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q2"
$newTotal.Range("C5").Value = 2
$newTotal.Range("D5").Value = 0.01

$newTotal.Range("A6").Value = 4
$newTotal.Range("B6").Value = "2021-Q1"
$newTotal.Range("C6").Value = 2
$newTotal.Range("D6").Value = 0.01
